# Update cryptocurrency price/volume data (cryptos list refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.839.95"
$ws.Range("D3").Value = "1.634.33"
$ws.Range("E3").Value = "  -1.46%  "
$ws.Range("E4").Value = "  -0.16%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "215.02"
$ws.Range("E5").Value = "  -1.39%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5012"
$ws.Range("E6").Value = "  -2.77%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2567"
$ws.Range("E8").Value = "  -0.69%  "
$ws.Range("E9").Value = "  -0.51%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.54"
$ws.Range("E10").Value = "  -2.13%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07680"
$ws.Range("E11").Value = "  -1.40%  "
$ws.Range("B12").Value = "Polkadot"
$ws.Range("C12").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D12").Value = "4.226"
$ws.Range("E12").Value = "  -1.65%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.626.66"
$ws.Range("E13").Value = "  -1.93%  "
$ws.Range("D14").Value = "1.859.65"
$ws.Range("E14").Value = "  -1.46%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5450"
$ws.Range("E15").Value = "  -1.80%  "
$ws.Range("D16").Value = "0.0₅7915"
$ws.Range("E16").Value = "  -1.95%  "
$ws.Range("D17").Value = "63.43"
$ws.Range("E17").Value = "  -1.37%  "
$ws.Range("D18").Value = "25.857.12"
$ws.Range("E18").Value = "  -1.41%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "202.95"
$ws.Range("E20").Value = "  -4.18%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.297"
$ws.Range("E21").Value = "  -2.92%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.929"
$ws.Range("E22").Value = "  -1.13%  "
$ws.Range("D23").Value = "5.968"
$ws.Range("E23").Value = "  +0.17%  "
$ws.Range("E24").Value = "  -0.11%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.936"
$ws.Range("E25").Value = "  +10.43%  "
$ws.Range("D26").Value = "140.88"
$ws.Range("E26").Value = "  -1.99%  "
$ws.Range("D27").Value = "0.1143"
$ws.Range("E27").Value = "  -1.78%  "
$ws.Range("D28").Value = "15.66"
$ws.Range("E28").Value = "  -0.92%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.692"
$ws.Range("E29").Value = "  -4.04%  "
$ws.Range("E30").Value = "  -1.31%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.04966"
$ws.Range("E31").Value = "  -5.49%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.273"
$ws.Range("E32").Value = "  -2.84%  "
$ws.Range("E33").Value = "  -1.34%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.527"
$ws.Range("E34").Value = "  -2.73%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.354"
$ws.Range("E35").Value = "  -0.63%  "
$ws.Range("D36").Value = "1.170.91"
$ws.Range("E36").Value = "  +0.48%  "
$ws.Range("E37").Value = "  -4.96%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.8903"
$ws.Range("E38").Value = "  -4.23%  "
$ws.Range("E39").Value = "  -2.21%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01553"
$ws.Range("E40").Value = "  -2.48%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.553"
$ws.Range("E41").Value = "  -0.39%  "
$ws.Range("E42").Value = "  -0.19%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.627"
$ws.Range("E43").Value = "  -0.98%  "
$ws.Range("E44").Value = "  -5.54%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "99.10"
$ws.Range("E45").Value = "  -1.27%  "
$ws.Range("D46").Value = "1.772.15"
$ws.Range("E46").Value = "  -1.38%  "
$ws.Range("D47").Value = "0.0₈111"
$ws.Range("E47").Value = "  -1.31%  "
$ws.Range("D48").Value = "0.4511"
$ws.Range("E48").Value = "  -0.52%  "
$ws.Range("D49").Value = "1.004"
$ws.Range("E49").Value = "  -0.08%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "54.71"
$ws.Range("E50").Value = "  -2.20%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05032"
$ws.Range("E51").Value = "  -0.42%  "
